$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("事業投資")

# --- Row 1: turn the old sample-data row into a proper header row ---
$ws.Range("B1").Value = "owner"
$ws.Range("C1").Value = "company"
$ws.Range("D1").Value = "address"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Give the newly added header cells (H1:N1) the same look as the
# pre-existing header cells (bold, centered, bordered).
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# --- Row 2: extend the data row with the common metadata columns ---
$ws.Range("H2").Value = "investment"
$ws.Range("I2").Value = "normal"

# Force J2 to stay plain text ("2012-04-20") instead of being
# auto-converted to a date serial number.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-20"
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = "廖正井"
$ws.Range("L2").Value = 1711
$ws.Range("M2").Value = "tmp845a1"
$ws.Range("N2").Value = 134
